# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Poroto verde" (Comercializadora del
# Agro de Limarí) at the top of that sub-block (rows 118-119), pushing the
# existing rows 118-146 down to 120-148.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 118 (shifts rows 118:146 -> 120:148)
$ws.Rows("118:119").Insert()

# New row 118: Magnum
$ws.Cells.Item(118, 1).Value2  = 2
$ws.Cells.Item(118, 2).Value2  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(118, 3).Value2  = "Coquimbo"
$ws.Cells.Item(118, 4).Value2  = 44588
$ws.Cells.Item(118, 5).Value2  = 4
$ws.Cells.Item(118, 6).Value2  = 100112031
$ws.Cells.Item(118, 7).Value2  = "Poroto verde"
$ws.Cells.Item(118, 8).Value2  = "Magnum"
$ws.Cells.Item(118, 9).Value2  = "Primera"
$ws.Cells.Item(118, 10).Value2 = 700
$ws.Cells.Item(118, 11).Value2 = 23000
$ws.Cells.Item(118, 12).Value2 = 24000
$ws.Cells.Item(118, 13).Value2 = 23500
$ws.Cells.Item(118, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(118, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(118, 16).Value2 = 940
$ws.Cells.Item(118, 17).Value2 = 25
$ws.Cells.Item(118, 18).Value2 = "Hortaliza"

# New row 119: Sin especificar
$ws.Cells.Item(119, 1).Value2  = 2
$ws.Cells.Item(119, 2).Value2  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(119, 3).Value2  = "Coquimbo"
$ws.Cells.Item(119, 4).Value2  = 44588
$ws.Cells.Item(119, 5).Value2  = 4
$ws.Cells.Item(119, 6).Value2  = 100112031
$ws.Cells.Item(119, 7).Value2  = "Poroto verde"
$ws.Cells.Item(119, 8).Value2  = "Sin especificar"
$ws.Cells.Item(119, 9).Value2  = "Primera"
$ws.Cells.Item(119, 10).Value2 = 708
$ws.Cells.Item(119, 11).Value2 = 26000
$ws.Cells.Item(119, 12).Value2 = 27000
$ws.Cells.Item(119, 13).Value2 = 26506
$ws.Cells.Item(119, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(119, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(119, 16).Value2 = 1060
$ws.Cells.Item(119, 17).Value2 = 25
$ws.Cells.Item(119, 18).Value2 = "Hortaliza"
